$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A5').Value = 'CM-6 b,SC-5,SC-5 (2)'
$ws.Range('A6').Value = 'AU-12 (3),AU-7 b,AC-6 (9),AU-8 b,CM-5 (1),AC-6 (8),AU-7 a'
$ws.Range('C6').Value = 'SRG-OS-000326-GPOS-00126,SRG-OS-000327-GPOS-00127,SRG-OS-000337-GPOS-00129,SRG-OS-000348-GPOS-00136,SRG-OS-000349-GPOS-00137,SRG-OS-000350-GPOS-00138,SRG-OS-000351-GPOS-00139,SRG-OS-000352-GPOS-00140,SRG-OS-000353-GPOS-00141,SRG-OS-000354-GPOS-00142,SRG-OS-000358-GPOS-00145,SRG-OS-000365-GPOS-00152'
$ws.Range('A7').Value = 'AU-12 (3),CM-6 b,AU-7 b,AU-8 b,CM-5 (1),AU-12 a,AU-7 a,AU-12 c'
$ws.Range('A11').Value = 'IA-2 (11),IA-2 (12)'
$ws.Range('A12').Value = 'IA-2 (11),IA-2 (12)'
$ws.Range('A17').Value = 'CM-6 b,CM-7 (2)'
$ws.Range('A22').Value = 'CM-6 b,CM-7 (2)'
$ws.Range('A23').Value = 'CM-6 b,CM-7 (2)'
$ws.Range('A38').Value = 'AC-7 a,AC-7 b'
$ws.Range('A39').Value = 'AC-7 a,AC-7 b'
$ws.Range('A40').Value = 'AC-7 a,AC-7 b'
$ws.Range('A41').Value = 'AC-7 a,AC-7 b'
$ws.Range('A45').Value = 'IA-8,AU-3 (1),IA-2'
$ws.Range('A46').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A47').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A48').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A49').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A50').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A51').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A52').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A53').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A54').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A55').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A56').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A57').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A58').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A59').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A60').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A61').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A62').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A63').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A64').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A65').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A66').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A67').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A68').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A69').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A70').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A71').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A72').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A73').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A74').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A75').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A76').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A77').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A78').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A79').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A80').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A81').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A82').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A83').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A84').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A85').Value = 'AU-3 (1),AU-3,MA-4 (1) (a)'
$ws.Range('A86').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A87').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A88').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A89').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A90').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A91').Value = 'AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A92').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A93').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A94').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A95').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A96').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A97').Value = 'AU-3 (1),AU-3,MA-4 (1) (a),AU-12 c'
$ws.Range('A98').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A99').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A100').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A101').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A102').Value = 'AU-3,AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A103').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A104').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A105').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A106').Value = 'AU-3,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A107').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A108').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A109').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A110').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A111').Value = 'AU-3,AU-12 a,AU-3 (1),AC-2 (4),MA-4 (1) (a),AU-12 c'
$ws.Range('A112').Value = 'AU-3,AU-14 (1),AU-12 a,AU-3 (1),MA-4 (1) (a),AU-12 c'
$ws.Range('A115').Value = 'AC-6 (10),AC-11 b'
$ws.Range('A120').Value = 'AU-3,MA-4 (1) (a),AU-12 a,AU-12 c'
$ws.Range('A121').Value = 'AU-3,MA-4 (1) (a),AU-12 a,AU-12 c'
$ws.Range('A122').Value = 'AU-9,AU-12 c'
$ws.Range('A127').Value = 'CM-5 (1),AC-2 (4),AC-6 (9),AU-12 c'
$ws.Range('A129').Value = 'CM-6 b,IA-5 (1) (a),IA-5 (1) (b)'
$ws.Range('A133').Value = 'AC-17 (2),SC-13,MA-4 c,SC-8'
$ws.Range('A134').Value = 'MA-4 (7),SC-10,MA-4 e,AC-12'
$ws.Range('A135').Value = 'SC-10,AC-12'
$ws.Range('A136').Value = 'SC-10,AC-12'
$ws.Range('A138').Value = 'AU-3,AU-6 (4),CM-6 b,AU-14 (1),CM-5 (1),AU-12 a,AU-3 (1),AU-7 a,MA-4 (1) (a),AU-7 (1)'
$ws.Range('A143').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A144').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A145').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A146').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A147').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A148').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A149').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A150').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A151').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A152').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A153').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A154').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A155').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A166').Value = 'SC-8 (1),SC-8 (2),SC-8'
$ws.Range('A169').Value = 'SC-8 (2),SC-8'
$ws.Range('A179').Value = 'CM-6 b,AU-6 (4),AU-4 (1)'
$ws.Range('A180').Value = 'CM-7 b,AC-17 (9),AC-17 (1),CM-6 b'
$ws.Range('A181').Value = 'CM-7 b,CM-6 b,AC-17 (1)'
$ws.Range('A200').Value = 'AU-3,AU-4 (1)'
$ws.Range('A207').Value = 'CM-6 b,AU-4 (1)'
$ws.Range('A208').Value = 'SC-28 (1),SC-28'
$ws.Range('A221').Value = 'CM-6 b,IA-2 (5)'
$ws.Range('A222').Value = 'IA-2,IA-2 (3),IA-2 (4),IA-2 (2),IA-2 (5)'
$ws.Range('A223').Value = 'IA-2,IA-2 (3),IA-2 (4),IA-2 (2),IA-2 (5)'
$ws.Range('A224').Value = 'AC-18 (1),SC-8,SC-8 (1)'
$ws.Range('A226').Value = 'IA-7,IA-5 (1) (c)'
$ws.Range('A227').Value = 'CM-6 b,IA-7'
$ws.Range('A228').Value = 'CM-6 b,IA-7'
$ws.Range('A229').Value = 'CM-6 b,IA-7'
$ws.Range('A231').Value = 'CM-7 a,IA-7'
$ws.Range('A233').Value = 'AC-17 (2),MA-4 (6)'
$ws.Range('A242').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A243').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A245').Value = 'CM-6 b,SC-2,SI-16'
$ws.Range('A247').Value = 'SC-3,SI-16'
$ws.Range('A259').Value = 'CM-6 b,IA-3'
$ws.Range('A260').Value = 'CM-6 b,IA-3'
$ws.Range('A261').Value = 'CM-6 b,IA-3'
$ws.Range('A262').Value = 'CM-6 b,IA-3'
$ws.Range('A270').Value = 'CM-6 b,IA-2 (2)'
$ws.Range('A271').Value = 'IA-2 (1),IA-2 (3),IA-2 (4),IA-2 (2)'
$ws.Range('A273').Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range('A277').Value = 'SC-4,SC-2'
$ws.Range('A278').Value = 'SC-4,SC-2'
$ws.Range('A281').Value = 'CM-6 b,AU-12 a'
$ws.Range('A298').Value = 'IA-2 (11),IA-2 (12)'
$ws.Range('A300').Value = 'IA-2 (1),IA-2 (11),IA-2 (12)'
$ws.Range('A310').Value = 'AU-8 b,AU-8 (1) (b),AU-8 (1) (a)'
$ws.Range('A322').Value = 'AC-3 (4),IA-11'
$ws.Range('A345').Value = 'CM-7 b,AC-17 (1)'
$ws.Range('A347').Value = 'CM-7 a,CM-6 b,IA-5 (1) (c)'
$ws.Range('A358').Value = 'AC-11 b,AC-11 (1)'
$ws.Range('A361').Value = 'SI-6 d,SI-6 b,CM-3 (5)'
$ws.Range('A367').Value = 'CM-7 a,SI-16'
$ws.Range('A385').Value = 'CM-6 b,AC-17 (2)'
$ws.Range('A391').Value = 'CM-6 b,IA-5 (1) (a)'
$ws.Range('A398').Value = 'CM-6 b,SI-16'
$ws.Range('A448').Value = 'CM-6 b,IA-5 (1) (c)'
$ws.Range('A524').Value = 'CM-6 b,SC-2'
$ws.Range('A525').Value = 'CM-6 b,SC-2'

$ws.Range('K7').Value = 'Run the following command to determine if the  audit  package is installed:  $ rpm -q audit 

If the audit package is not installed then this is a finding.'